# Updated cryptos list (Price / Volume(1h) refresh, plus the
# EnergySwap/WhiteBITCoin row re-ranking at rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds free-form text (e.g. "60.139.84", subscript
# notation like "0.0₃0791") -- force Text format so Excel does not
# reinterpret/round it as a number and drop significant digits.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "60.470.30"
$ws.Cells.Item(2, 5).Value = "  +3.51%  "
$ws.Cells.Item(3, 4).Value = "2.623.41"
$ws.Cells.Item(3, 5).Value = "  +1.67%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "
$ws.Cells.Item(5, 4).Value = "571.11"
$ws.Cells.Item(5, 5).Value = "  +6.22%  "
$ws.Cells.Item(6, 4).Value = "145.93"
$ws.Cells.Item(6, 5).Value = "  +1.54%  "
$ws.Cells.Item(7, 4).Value = "0.995"
$ws.Cells.Item(7, 5).Value = "  -0.43%  "
$ws.Cells.Item(8, 4).Value = "0.601"
$ws.Cells.Item(8, 5).Value = "  +4.11%  "
$ws.Cells.Item(9, 4).Value = "2.644.29"
$ws.Cells.Item(9, 5).Value = "  +2.38%  "
$ws.Cells.Item(10, 4).Value = "6.76"
$ws.Cells.Item(10, 5).Value = "  -0.57%  "
$ws.Cells.Item(11, 4).Value = "0.105"
$ws.Cells.Item(11, 5).Value = "  +4.40%  "
$ws.Cells.Item(12, 4).Value = "0.152"
$ws.Cells.Item(12, 5).Value = "  +10.25%  "
$ws.Cells.Item(13, 4).Value = "0.344"
$ws.Cells.Item(13, 5).Value = "  +3.47%  "
$ws.Cells.Item(14, 4).Value = "3.082.39"
$ws.Cells.Item(14, 5).Value = "  +1.36%  "
$ws.Cells.Item(15, 4).Value = "60.383.51"
$ws.Cells.Item(15, 5).Value = "  +3.42%  "
$ws.Cells.Item(16, 4).Value = "22.13"
$ws.Cells.Item(16, 5).Value = "  +6.90%  "
$ws.Cells.Item(17, 4).Value = "0.0000138"
$ws.Cells.Item(17, 5).Value = "  +3.86%  "
$ws.Cells.Item(18, 4).Value = "2.633.59"
$ws.Cells.Item(18, 5).Value = "  +1.31%  "
$ws.Cells.Item(19, 4).Value = "4.54"
$ws.Cells.Item(19, 5).Value = "  +1.41%  "
$ws.Cells.Item(20, 4).Value = "341.99"
$ws.Cells.Item(20, 5).Value = "  +1.71%  "
$ws.Cells.Item(21, 4).Value = "10.41"
$ws.Cells.Item(21, 5).Value = "  +3.56%  "
$ws.Cells.Item(22, 4).Value = "6.37"
$ws.Cells.Item(22, 5).Value = "  +3.38%  "
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).Value = "65.83"
$ws.Cells.Item(24, 5).Value = "  -1.57%  "
$ws.Cells.Item(25, 4).Value = "0.450"
$ws.Cells.Item(25, 5).Value = "  +7.27%  "
$ws.Cells.Item(26, 5).Value = "  +3.36%  "
$ws.Cells.Item(27, 5).Value = "  -0.14%  "
$ws.Cells.Item(28, 4).Value = "7.37"
$ws.Cells.Item(28, 5).Value = "  +4.50%  "
$ws.Cells.Item(29, 4).Value = "0.0₃0798"
$ws.Cells.Item(29, 5).Value = "  +8.43%  "
$ws.Cells.Item(30, 4).Value = "0.997"
$ws.Cells.Item(30, 5).Value = "  -0.19%  "
$ws.Cells.Item(31, 5).Value = "  +4.29%  "
$ws.Cells.Item(32, 4).Value = "6.15"
$ws.Cells.Item(32, 5).Value = "  +3.16%  "
$ws.Cells.Item(33, 4).Value = "159.63"
$ws.Cells.Item(33, 5).Value = "  +2.34%  "
$ws.Cells.Item(34, 4).Value = "19.18"
$ws.Cells.Item(34, 5).Value = "  +1.77%  "
$ws.Cells.Item(35, 4).Value = "4.11"
$ws.Cells.Item(35, 5).Value = "  +5.53%  "
$ws.Cells.Item(36, 4).Value = "1.15"
$ws.Cells.Item(36, 5).Value = "  +4.74%  "
$ws.Cells.Item(37, 4).Value = "0.891"
$ws.Cells.Item(37, 5).Value = "  +8.88%  "
$ws.Cells.Item(38, 4).Value = "0.885"
$ws.Cells.Item(38, 5).Value = "  +4.73%  "
$ws.Cells.Item(39, 4).Value = "37.56"
$ws.Cells.Item(39, 5).Value = "  +1.64%  "
$ws.Cells.Item(40, 4).Value = "1.51"
$ws.Cells.Item(40, 5).Value = "  +6.97%  "
$ws.Cells.Item(41, 4).Value = "297.76"
$ws.Cells.Item(41, 5).Value = "  +6.43%  "
$ws.Cells.Item(42, 4).Value = "3.66"
$ws.Cells.Item(42, 5).Value = "  +1.85%  "
$ws.Cells.Item(43, 4).Value = "0.995"
$ws.Cells.Item(43, 5).Value = "  -0.36%  "
$ws.Cells.Item(44, 4).Value = "0.0982"
$ws.Cells.Item(44, 5).Value = "  +4.61%  "
$ws.Cells.Item(45, 4).Value = "0.603"
$ws.Cells.Item(45, 5).Value = "  +2.19%  "
$ws.Cells.Item(46, 4).Value = "0.0542"
$ws.Cells.Item(46, 5).Value = "  +1.61%  "
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value = "19.33"
$ws.Cells.Item(47, 5).Value = "  +4.47%  "
$ws.Cells.Item(48, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(48, 4).Value = "10.67"
$ws.Cells.Item(48, 5).Value = "  +0.24%  "
$ws.Cells.Item(49, 4).Value = "126.85"
$ws.Cells.Item(49, 5).Value = "  +16.03%  "
$ws.Cells.Item(50, 4).Value = "0.0237"
$ws.Cells.Item(50, 5).Value = "  +4.31%  "
$ws.Cells.Item(51, 4).Value = "4.67"
$ws.Cells.Item(51, 5).Value = "  +6.15%  "
